$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Jan_2026), shifting D:H -> E:I
$ws.Range("D1").EntireColumn.Insert()

# Copy the header formatting from the neighboring column (C1) onto the new D1 cell
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new column header
$ws.Range("D1").Value = "Status"

# The old "Nov_2025" header (now shifted to G1) becomes "Oct_2025"
$ws.Range("G1").Value = "Oct_2025"
